$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column A, shifting existing columns (country, gini,
# hdi, gdpcap, corruption) one to the right.
[void]$ws.Columns.Item(1).EntireColumn.Insert()

# Header for the newly inserted column A.
$ws.Range("A1").Value = "Code"

# Fill in the 3-letter country-code lookup formula for every data row, now that
# the country name lives in column B instead of column A. Row 49 (Russian
# Federation) is handled separately below.
for ($i = 2; $i -le 63; $i++) {
    if ($i -eq 49) {
        continue
    }
    $ws.Range("A$i").Formula = "=VLOOKUP(B$i,[1]Sheet2!`$B`$5:`$C`$181,2,FALSE)"
}

# The external lookup table lists this country as "Russia" rather than
# "Russian Federation", so the VLOOKUP can't resolve it there; the code was
# entered as a literal value instead of a formula.
$ws.Range("A49").Value = "RUS"

# Re-apply the autofilter so it covers the new column (A1:F63 instead of A1:E63).
if ($ws.AutoFilterMode) {
    $ws.AutoFilterMode = $false
}
[void]$ws.Range("A1:F63").AutoFilter()

# Keep the hidden _FilterDatabase defined name in sync with the wider autofilter range.
$filterName = $wb.Names.Item("_xlnm._FilterDatabase")
$filterName.RefersTo = "=pyrun!`$A`$1:`$F`$63"

# Update the active selection to match the saved view.
[void]$ws.Range("B11").Select()
